# aggiornamento fino a 02/05
# Append 6 new rows (239-244) of data to the bottom of the table, continuing
# the daily series that ends at row 238 (date serial 44312 / 2021-04-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date serial, nuovi pos. (B), somma mobile 7gg. (C),
# somma mobile 7gg. per 100mila abitanti (D)
$rows = @(
    @{ R = 239; A = 44313; B = 0; C = 9;  D = 274.8091603053435 },
    @{ R = 240; A = 44314; B = 1; C = 10; D = 305.3435114503817 },
    @{ R = 241; A = 44315; B = 7; C = 15; D = 458.0152671755725 },
    @{ R = 242; A = 44316; B = 5; C = 19; D = 580.1526717557252 },
    @{ R = 243; A = 44317; B = 3; C = 21; D = 641.2213740458016 },
    @{ R = 244; A = 44318; B = 1; C = 20; D = 610.6870229007634 }
)

foreach ($row in $rows) {
    $r = $row.R

    # Column A carries the date style (s="2") used by every row above it;
    # copy formatting from the prior row's A cell, then overwrite the value.
    $ws.Cells.Item($r - 1, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row.A

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
}
